$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = (Get-Date -Year 2024 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0)

# --- "MCF" sheet: set all non-zero capacity factors (other than nuclear, which is
#     already 1) to 1 ---
$wsMCF = $wb.Worksheets.Item("MCF")

$wsMCF.Range("B2").Value = 1
$wsMCF.Range("B3").Value = 1
$wsMCF.Range("B4").Value = 1
$wsMCF.Range("B6").Value = 1
$wsMCF.Range("B10").Value = 1
$wsMCF.Range("B11").Value = 1
$wsMCF.Range("B12").Value = 1
$wsMCF.Range("B13").Value = 1
$wsMCF.Range("B14").Value = 1
$wsMCF.Range("B16").Value = 1
$wsMCF.Range("B17").Value = 1
$wsMCF.Range("B18").Value = 1

# Update the selected cell on the MCF sheet to match the saved view state
$wsMCF.Activate()
$wsMCF.Range("B17").Select()

$wb.Save()
